$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("global_settings")

# Enable FTUX settings: mapAsButton (F11), showNextDragonInXpBar (D11), showUnlockProgressionText (E11)
$ws.Range("D11").Value = $true
$ws.Range("E11").Value = $true
$ws.Range("F11").Value = $true

# Update the saved selection on the sheet view to D14
$ws.Range("D14").Select()
